$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell. Values that look like a plain
# decimal number (single dot, optional sign, no other characters) would
# otherwise be auto-converted to a numeric type by Excel (losing trailing
# zeros / thousands-style dots), so those are forced to Text format first.
function Set-Text($range, $value) {
    $r = $ws.Range($range)
    if ($value -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $r.NumberFormat = "@"
    }
    $r.Value = $value
}

# Row 2 - Bitcoin
Set-Text "D2" "26.973.22"
Set-Text "E2" "  -0.94%  "

# Row 3 - Ethereum
Set-Text "D3" "1.811.42"
Set-Text "E3" "  -0.45%  "

# Row 4 - TetherUSD
Set-Text "D4" "1.001"
Set-Text "E4" "  -0.02%  "

# Row 5 - BNB
Set-Text "D5" "310.68"
Set-Text "E5" "  -0.77%  "

# Row 6 - USDC (D unchanged)
Set-Text "E6" "  +0.01%  "

# Row 7 - XRP
Set-Text "D7" "0.4629"
Set-Text "E7" "  +3.74%  "

# Row 8 - Cardano
Set-Text "D8" "0.3720"
Set-Text "E8" "  -1.04%  "

# Row 9 - Dogecoin
Set-Text "D9" "0.07380"
Set-Text "E9" "  -0.03%  "

# Row 10 - Polygon
Set-Text "D10" "0.8757"
Set-Text "E10" "  -0.31%  "

# Row 11 - Solana
Set-Text "D11" "20.49"
Set-Text "E11" "  -1.56%  "

# Row 12 - WrappedEther
Set-Text "D12" "1.814.81"
Set-Text "E12" "  -0.27%  "

# Row 13 - Polkadot
Set-Text "D13" "5.367"
Set-Text "E13" "  -1.00%  "

# Row 14 - Litecoin
Set-Text "D14" "92.59"
Set-Text "E14" "  -0.25%  "

# Row 15 - Chainlink
Set-Text "D15" "6.536"
Set-Text "E15" "  -2.36%  "

# Row 16 - TRON
Set-Text "D16" "0.07065"
Set-Text "E16" "  -0.18%  "

# Row 17 - BinanceUSD
Set-Text "D17" "1.002"
Set-Text "E17" "  -0.01%  "

# Row 18 - ShibaInu
Set-Text "D18" "0.000008718"
Set-Text "E18" "  -0.86%  "

# Row 19 - Dai
Set-Text "D19" "1.002"
Set-Text "E19" "  +0.10%  "

# Row 20 - Avalanche
Set-Text "D20" "14.75"
Set-Text "E20" "  -1.65%  "

# Row 21 - WrappedBTC
Set-Text "D21" "26.965.12"
Set-Text "E21" "  -1.03%  "

# Row 22 - Uniswap
Set-Text "D22" "5.317"
Set-Text "E22" "  -0.42%  "

# Row 23 - Cosmos (D unchanged)
Set-Text "E23" "  -2.42%  "

# Row 24 - WrappedliquidstakedEther2.0
Set-Text "D24" "2.035.24"
Set-Text "E24" "  -0.44%  "

# Row 25 - Toncoin
Set-Text "D25" "1.898"
Set-Text "E25" "  -3.17%  "

# Row 26 - Monero
Set-Text "D26" "151.60"
Set-Text "E26" "  +0.40%  "

# Row 27 - EthereumClassic
Set-Text "D27" "18.35"
Set-Text "E27" "  -1.19%  "

# Row 28 - LidoDAOToken
Set-Text "D28" "2.154"
Set-Text "E28" "  -5.70%  "

# Row 29 - InternetComputer(DFINITY)
Set-Text "D29" "5.289"
Set-Text "E29" "  -1.02%  "

# Row 30 - BitcoinCash
Set-Text "D30" "115.95"
Set-Text "E30" "  -0.97%  "

# Row 31 - Stellar
Set-Text "D31" "0.08936"
Set-Text "E31" "  +0.87%  "

# Row 32 - ImmutableX
Set-Text "D32" "0.7556"
Set-Text "E32" "  -4.25%  "

# Row 33 - ARBITRUM
Set-Text "D33" "1.158"
Set-Text "E33" "  -2.78%  "

# Row 34 - HuobiToken
Set-Text "D34" "2.925"
Set-Text "E34" "  +0.02%  "

# Row 35 - Filecoin
Set-Text "D35" "4.460"
Set-Text "E35" "  -2.19%  "

# Row 36 - Frax (D unchanged)
Set-Text "E36" "  +0.03%  "

# Row 37 - TrustWalletToken
Set-Text "D37" "1.107"
Set-Text "E37" "  -0.02%  "

# Row 38 - VeChain
Set-Text "D38" "0.01978"
Set-Text "E38" "  +0.29%  "

# Row 39 - RenderToken
Set-Text "D39" "2.460"
Set-Text "E39" "  +6.44%  "

# Row 40 - Hedera
Set-Text "D40" "0.05247"
Set-Text "E40" "  -0.17%  "

# Row 41 - MXToken
Set-Text "D41" "2.922"
Set-Text "E41" "  +1.83%  "

# Row 42 - TheSandbox
Set-Text "D42" "0.5322"
Set-Text "E42" "  +0.83%  "

# Row 43 - FraxShare
Set-Text "D43" "7.203"
Set-Text "E43" "  -1.15%  "

# Row 44 - Algorand
Set-Text "D44" "0.1664"
Set-Text "E44" "  -2.06%  "

# Row 45 - Aptos
Set-Text "D45" "8.535"
Set-Text "E45" "  -0.79%  "

# Row 46 - Decentraland
Set-Text "D46" "0.4993"
Set-Text "E46" "  -0.83%  "

# Row 47 - EnergySwap
Set-Text "D47" "10.36"
Set-Text "E47" "  -1.65%  "

# Row 48 - Quant
Set-Text "D48" "104.30"
Set-Text "E48" "  -0.51%  "

# Row 49 - was PaxDollar, now NEARProtocol (rows 49/50 swap content)
Set-Text "B49" "NEARProtocol"
Set-Text "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-Text "D49" "1.675"
Set-Text "E49" "  -0.60%  "

# Row 50 - was NEARProtocol, now PaxDollar
Set-Text "B50" "PaxDollar"
Set-Text "C50" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-Text "D50" "1.000"
Set-Text "E50" "  +0.02%  "

# Row 51 - Cronos
Set-Text "D51" "0.06298"
Set-Text "E51" "  -1.30%  "
